$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.054610263093101
$ws.Range("D2").Value = 1.053306798612993
$ws.Range("E2").Value = 1.068088514750336
$ws.Range("F2").Value = 1.07512130141437
$ws.Range("I2").Value = 1.04859352851958
$ws.Range("J2").Value = 1.059621457931861
$ws.Range("K2").Value = 1.05605312985389
$ws.Range("L2").Value = 1.070794633638798
$ws.Range("M2").Value = 1.077808707890674
$ws.Range("N2").Value = 1.023525571709366

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.055802924411355
$ws.Range("D3").Value = 1.054209404730506
$ws.Range("E3").Value = 1.069305707358738
$ws.Range("F3").Value = 1.076480496108739
$ws.Range("I3").Value = 1.048983139922067
$ws.Range("J3").Value = 1.060464291693418
$ws.Range("K3").Value = 1.056768524974702
$ws.Range("L3").Value = 1.07182673942624
$ws.Range("M3").Value = 1.078983829075165
$ws.Range("N3").Value = 1.023817917550371

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056574196938277
$ws.Range("D4").Value = 1.054792966454642
$ws.Range("E4").Value = 1.070093241359226
$ws.Range("F4").Value = 1.077360188796475
$ws.Range("I4").Value = 1.049233662641216
$ws.Range("J4").Value = 1.061008659025667
$ws.Range("K4").Value = 1.057230316573381
$ws.Range("L4").Value = 1.072493938871682
$ws.Range("M4").Value = 1.079743869130526
$ws.Range("N4").Value = 1.024006476319868

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.056898332245845
$ws.Range("D5").Value = 1.055038180739469
$ws.Range("E5").Value = 1.070424305615185
$ws.Range("F5").Value = 1.077730062560424
$ws.Range("I5").Value = 1.049338604416671
$ws.Range("J5").Value = 1.061237272487854
$ws.Range("K5").Value = 1.05742418747344
$ws.Range("L5").Value = 1.072774277584035
$ws.Range("M5").Value = 1.08006331019014
$ws.Range("N5").Value = 1.024085601109145

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056952749742106
$ws.Range("D6").Value = 1.05507934657834
$ws.Range("E6").Value = 1.070479891935669
$ws.Range("F6").Value = 1.07779216905083
$ws.Range("I6").Value = 1.049356202459256
$ws.Range("J6").Value = 1.061275643746071
$ws.Range("K6").Value = 1.05745672366717
$ws.Range("L6").Value = 1.072821338829137
$ws.Range("M6").Value = 1.080116941131733
$ws.Range("N6").Value = 1.024098877998627

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056578528470351
$ws.Range("D7").Value = 1.054796243470801
$ws.Range("E7").Value = 1.070097665112651
$ws.Range("F7").Value = 1.077365130864578
$ws.Range("I7").Value = 1.049235066362624
$ws.Range("J7").Value = 1.061011714705549
$ws.Range("K7").Value = 1.057232908130802
$ws.Range("L7").Value = 1.072497685366071
$ws.Range("M7").Value = 1.079748137828479
$ws.Range("N7").Value = 1.024007534159711

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.055013424981424
$ws.Range("D8").Value = 1.053611939556155
$ws.Range("E8").Value = 1.068499886575798
$ws.Range("F8").Value = 1.075580607071869
$ws.Range("I8").Value = 1.048725527506513
$ws.Range("J8").Value = 1.059906505976993
$ws.Range("K8").Value = 1.056295132902799
$ws.Range("L8").Value = 1.071143572618507
$ws.Range("M8").Value = 1.078205917647692
$ws.Range("N8").Value = 1.023624497578776

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.052251895681693
$ws.Range("D9").Value = 1.051521286688022
$ws.Range("E9").Value = 1.065683764608641
$ws.Range("F9").Value = 1.072437493912655
$ws.Range("I9").Value = 1.047815506133531
$ws.Range("J9").Value = 1.057951244912742
$ws.Range("K9").Value = 1.054634049691265
$ws.Range("L9").Value = 1.068752456512417
$ws.Range("M9").Value = 1.075485616487785
$ws.Range("N9").Value = 1.022944860852314

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.050408301085925
$ws.Range("D10").Value = 1.05012492569641
$ws.Range("E10").Value = 1.063805799281021
$ws.Range("F10").Value = 1.070342904301423
$ws.Range("I10").Value = 1.047200609170978
$ws.Range("J10").Value = 1.056642440159757
$ws.Range("K10").Value = 1.053520805222448
$ws.Range("L10").Value = 1.067154906593047
$ws.Range("M10").Value = 1.073670126061205
$ws.Range("N10").Value = 1.022488600695608

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.049609360942193
$ws.Range("D11").Value = 1.049519655656095
$ws.Range("E11").Value = 1.062992460213593
$ws.Range("F11").Value = 1.069436083425846
$ws.Range("I11").Value = 1.046932392161675
$ws.Range("J11").Value = 1.056074435648743
$ws.Range("K11").Value = 1.053037352706644
$ws.Range("L11").Value = 1.066462300319329
$ws.Range("M11").Value = 1.072883504388871
$ws.Range("N11").Value = 1.02229027767013

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049312498216509
$ws.Range("D12").Value = 1.04929473435634
$ws.Range("E12").Value = 1.062690322105488
$ws.Range("F12").Value = 1.069099269263621
$ws.Range("I12").Value = 1.046832468541781
$ws.Range("J12").Value = 1.055863258752019
$ws.Range("K12").Value = 1.052857563523701
$ws.Range("L12").Value = 1.066204904489239
$ws.Range("M12").Value = 1.072591239903366
$ws.Range("N12").Value = 1.022216497049027

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049376180886221
$ws.Range("D13").Value = 1.049342985153705
$ws.Range("E13").Value = 1.062755133053095
$ws.Range("F13").Value = 1.069171516212576
$ws.Range("I13").Value = 1.046853915903176
$ws.Range("J13").Value = 1.055908565755984
$ws.Range("K13").Value = 1.052896138594525
$ws.Range("L13").Value = 1.066260122734965
$ws.Range("M13").Value = 1.072653935208542
$ws.Range("N13").Value = 1.022232328438044

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049584824249437
$ws.Range("D14").Value = 1.049501065579516
$ws.Range("E14").Value = 1.062967485964615
$ws.Range("F14").Value = 1.069408241882474
$ws.Range("I14").Value = 1.0469241384859
$ws.Range("J14").Value = 1.056056983684116
$ws.Range("K14").Value = 1.053022495638595
$ws.Range("L14").Value = 1.066441026585995
$ws.Range("M14").Value = 1.072859347325666
$ws.Range("N14").Value = 1.022284181282257

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.049713362869928
$ws.Range("D15").Value = 1.049598451243311
$ws.Range("E15").Value = 1.063098319845802
$ws.Range("F15").Value = 1.069554098870409
$ws.Range("I15").Value = 1.046967365698978
$ws.Range("J15").Value = 1.056148403007001
$ws.Range("K15").Value = 1.053100320061888
$ws.Range("L15").Value = 1.066552469987875
$ws.Range("M15").Value = 1.072985898070976
$ws.Range("N15").Value = 1.022316114326819

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.050461310137921
$ws.Range("D16").Value = 1.050165081943231
$ws.Range("E16").Value = 1.063859774163988
$ws.Range("F16").Value = 1.070403089918756
$ws.Range("I16").Value = 1.04721836841394
$ws.Range("J16").Value = 1.056680109545312
$ws.Range("K16").Value = 1.053552860531342
$ws.Range("L16").Value = 1.067200854344206
$ws.Range("M16").Value = 1.073722320736846
$ws.Range("N16").Value = 1.022501746697233

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.050930301223261
$ws.Range("D17").Value = 1.050520342991716
$ws.Range("E17").Value = 1.064337367632382
$ws.Range("F17").Value = 1.070935677622618
$ws.Range("I17").Value = 1.04737528962662
$ws.Range("J17").Value = 1.05701329005401
$ws.Range("K17").Value = 1.053836348566812
$ws.Range("L17").Value = 1.06760733774086
$ws.Range("M17").Value = 1.074184122823864
$ws.Range("N17").Value = 1.022617985299133

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.051203793124053
$ws.Range("D18").Value = 1.0507274994341
$ws.Range("E18").Value = 1.064615923818101
$ws.Range("F18").Value = 1.071246341726499
$ws.Range("I18").Value = 1.04746662982594
$ws.Range("J18").Value = 1.057207504907152
$ws.Range("K18").Value = 1.054001566369886
$ws.Range("L18").Value = 1.067844350136383
$ws.Range("M18").Value = 1.074453435883908
$ws.Range("N18").Value = 1.02268571207764

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.051297036252718
$ws.Range("D19").Value = 1.050798124069827
$ws.Range("E19").Value = 1.064710901641498
$ws.Range("F19").Value = 1.071352272744913
$ws.Range("I19").Value = 1.047497742397527
$ws.Range("J19").Value = 1.057273706228712
$ws.Range("K19").Value = 1.054057878353162
$ws.Range("L19").Value = 1.067925151315923
$ws.Range("M19").Value = 1.074545256569015
$ws.Range("N19").Value = 1.022708792753511

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.050879989422184
$ws.Range("D20").Value = 1.050482233196847
$ws.Range("E20").Value = 1.064286128048291
$ws.Range("F20").Value = 1.070878534504788
$ws.Range("I20").Value = 1.047358473052613
$ws.Range("J20").Value = 1.056977555748298
$ws.Range("K20").Value = 1.053805947059418
$ws.Range("L20").Value = 1.067563734469801
$ws.Range("M20").Value = 1.074134580845373
$ws.Range("N20").Value = 1.022605521576068

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049523386784721
$ws.Range("D21").Value = 1.049454517524924
$ws.Range("E21").Value = 1.06290495411206
$ws.Range("F21").Value = 1.069338531559809
$ws.Range("I21").Value = 1.046903467861547
$ws.Range("J21").Value = 1.056013283694418
$ws.Range("K21").Value = 1.052985292541321
$ws.Range("L21").Value = 1.066387758546327
$ws.Range("M21").Value = 1.072798860756605
$ws.Range("N21").Value = 1.022268915080374

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.048669852318872
$ws.Range("D22").Value = 1.048807789636566
$ws.Range("E22").Value = 1.062036393323297
$ws.Range("F22").Value = 1.068370381797992
$ws.Range("I22").Value = 1.04661567565603
$ws.Range("J22").Value = 1.0554058799776
$ws.Range("K22").Value = 1.052468079097214
$ws.Range("L22").Value = 1.065647616600986
$ws.Range("M22").Value = 1.071958586899217
$ws.Range("N22").Value = 1.022056613800439

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049122383467656
$ws.Range("D23").Value = 1.049150686101113
$ws.Range("E23").Value = 1.062496850040109
$ws.Range("F23").Value = 1.068883606796271
$ws.Range("I23").Value = 1.046768402407352
$ws.Range("J23").Value = 1.055727983704733
$ws.Range("K23").Value = 1.052742381335124
$ws.Range("L23").Value = 1.066040052592854
$ws.Range("M23").Value = 1.072404075872304
$ws.Range("N23").Value = 1.022169221788571

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.050902723353038
$ws.Range("D24").Value = 1.050499453562261
$ws.Range("E24").Value = 1.064309281060217
$ws.Range("F24").Value = 1.070904354974574
$ws.Range("I24").Value = 1.047366072323427
$ws.Range("J24").Value = 1.056993702925106
$ws.Range("K24").Value = 1.053819684612147
$ws.Range("L24").Value = 1.067583437165999
$ws.Range("M24").Value = 1.074156966882119
$ws.Range("N24").Value = 1.022611153622325

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.052966260551665
$ws.Range("D25").Value = 1.052062222653874
$ws.Range("E25").Value = 1.066411887460381
$ws.Range("F25").Value = 1.073249908391103
$ws.Range("I25").Value = 1.048052213085966
$ws.Range("J25").Value = 1.058457653697076
$ws.Range("K25").Value = 1.055064506379253
$ws.Range("L25").Value = 1.069371221076101
$ws.Range("M25").Value = 1.076189214599738
$ws.Range("N25").Value = 1.023121120012596

Write-Output "vm_pu values updated"
